$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: some "Price" (column D) values are digit strings that use "."
# as a thousands separator (e.g. "68.873.65"), which Excel would
# otherwise misinterpret/auto-convert as a number when plausible
# (e.g. "559.21" -> 559.21). Prefixing the literal with an apostrophe
# forces Excel to store it as text (matching the source inlineStr
# cells), then ClearFormats() removes the quote-prefix style marker
# so the cell style is left exactly as it was (no explicit style).

$ws.Range("D2").Value = '68.873.65'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.440.28'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'559.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = "'162.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.11%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("E9").Value = '  +9.27%  '
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  -5.01%  '
$ws.Range("E13").Value = '  +4.65%  '
$ws.Range("D14").Value = '68.762.57'
$ws.Range("D15").Value = '2.888.40'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = '2.441.41'
$ws.Range("E17").Value = '  +3.14%  '
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").Value = "'339.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("E22").Value = '  +3.94%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = "'66.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = "'3.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '2.567.20'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = "'0.966"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.20%  '
$ws.Range("D29").Value = '0.0₃0823'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +2.46%  '
$ws.Range("D33").Value = "'430.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("D35").Value = "'159.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").Value = "'19.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = "'18.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").Value = "'4.38"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").Value = "'1.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.57%  '
$ws.Range("D44").Value = "'2.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").Value = "'130.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").Value = "'0.559"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("E50").Value = '  +2.97%  '
$ws.Range("E51").Value = '  +0.26%  '
